# functions.xlsx -- "Implemented H2 and HyperSQL dialects"
#
# The underlying dialect-comparison table already listed H2/HyperSQL; this
# pass brings the parameter-naming convention in the "Strings" section
# (instr/locate/strpos/charindex) in line with the rest of the sheet
# (renaming the optional "pos" parameter to "from"), and fixes the H2
# trunc() signature to match the others ("trunc(x, places)" -- no longer
# marked optional). It also nudges the window's zoom/scroll position to
# where the author was working (around the Strings block) when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Arithmetic section: H2's trunc() signature -------------------------
$ws.Range("E34").Value = "trunc(x, places)"

# --- Strings section: "pos" -> "from" for the optional search-start arg -
$ws.Range("F45").Value = "instr(t, sub [, from]) *1"       # Oracle
$ws.Range("F46").Value = "locate(sub, t [, from]) *1"      # DB2
$ws.Range("F47").Value = "strpos(t, from) *3"              # PostgreSQL
$ws.Range("F48").Value = "charindex(sub, t [, from]) *1"   # SQL Server
$ws.Range("F49").Value = "locate(sub, t [, from]) *1"      # MariaDB
$ws.Range("F50").Value = "locate(sub, t [, from]) *1"      # MySQL
$ws.Range("F51").Value = "charindex(sub, t [, from]) *1"   # SAP ASE 16
$ws.Range("F52").Value = "locate(sub, t [, from]) *1"      # H2
$ws.Range("F53").Value = "locate(sub, t [, from]) *1"      # HyperSQL
$ws.Range("F54").Value = "locate(sub, t [, from]) *1"      # Derby

# --- Minor row re-layout following the E34 edit --------------------------
$ws.Rows.Item(34).RowHeight = 12.95

# --- Window state: zoomed out slightly and scrolled to the edited area ---
$win = $excel.ActiveWindow
$win.Zoom = 110
$win.ScrollRow = 42
$win.ScrollColumn = 1
$ws.Range("A79").Select() | Out-Null

$win.TabRatio = 0.99
